# Updated cryptos list (GitHub Actions refresh): new Price/Volume(1h)
# snapshots for rows 2-51, plus row 51 now lists Decentraland instead of
# Algorand. Numeric-looking Price values are entered with a leading "'"
# so Excel keeps storing them as text (matching the original inline
# string cells) instead of auto-converting them to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.457.69"
$ws.Range("E2").Value = "  +0.59%  "

$ws.Range("D3").Value = "1.887.01"
$ws.Range("E3").Value = "  +1.01%  "

$ws.Range("D4").Value = "'0.9995"
$ws.Range("E4").Value = "  -0.11%  "

$ws.Range("D5").Value = "'245.08"
$ws.Range("E5").Value = "  +4.38%  "

$ws.Range("D6").Value = "'0.9994"
$ws.Range("E6").Value = "  -0.10%  "

$ws.Range("D7").Value = "'0.4739"
$ws.Range("E7").Value = "  +0.87%  "

$ws.Range("D8").Value = "'0.2907"
$ws.Range("E8").Value = "  +1.79%  "

$ws.Range("D9").Value = "'0.06541"
$ws.Range("E9").Value = "  -0.37%  "

$ws.Range("D10").Value = "'21.33"
$ws.Range("E10").Value = "  -0.02%  "

$ws.Range("D11").Value = "'0.07780"
$ws.Range("E11").Value = "  -0.49%  "

$ws.Range("D12").Value = "1.882.95"
$ws.Range("E12").Value = "  +0.59%  "

$ws.Range("D13").Value = "'0.7379"
$ws.Range("E13").Value = "  +6.44%  "

$ws.Range("D14").Value = "'96.16"
$ws.Range("E14").Value = "  -0.75%  "

$ws.Range("D15").Value = "'5.157"
$ws.Range("E15").Value = "  +1.53%  "

$ws.Range("D16").Value = "'275.78"
$ws.Range("E16").Value = "  +2.71%  "

$ws.Range("D17").Value = "30.448.26"
$ws.Range("E17").Value = "  +0.26%  "

$ws.Range("D18").Value = "'13.49"
$ws.Range("E18").Value = "  -1.94%  "

$ws.Range("D19").Value = "'0.000007600"
$ws.Range("E19").Value = "  -1.24%  "

$ws.Range("D20").Value = "'0.9994"
$ws.Range("E20").Value = "  -0.06%  "

$ws.Range("D21").Value = "2.129.78"
$ws.Range("E21").Value = "  +1.26%  "

$ws.Range("D22").Value = "'5.295"
$ws.Range("E22").Value = "  +0.76%  "

$ws.Range("D23").Value = "'0.9997"
$ws.Range("E23").Value = "  -0.08%  "

$ws.Range("D24").Value = "'6.205"
$ws.Range("E24").Value = "  +0.70%  "

$ws.Range("D25").Value = "'9.311"
$ws.Range("E25").Value = "  -2.68%  "

$ws.Range("D26").Value = "'165.08"
$ws.Range("E26").Value = "  -0.66%  "

$ws.Range("D27").Value = "'19.03"
$ws.Range("E27").Value = "  +0.78%  "

$ws.Range("D28").Value = "'1.963"
$ws.Range("E28").Value = "  +1.26%  "

$ws.Range("D29").Value = "'1.383"
$ws.Range("E29").Value = "  +1.43%  "

$ws.Range("D30").Value = "'0.09948"
$ws.Range("E30").Value = "  +0.51%  "

$ws.Range("D31").Value = "'1.520"
$ws.Range("E31").Value = "  +4.36%  "

$ws.Range("D32").Value = "'4.348"
$ws.Range("E32").Value = "  -0.07%  "

$ws.Range("E33").Value = "  +0.55%  "

$ws.Range("D34").Value = "'0.04786"
$ws.Range("E34").Value = "  +0.89%  "

$ws.Range("D35").Value = "'1.134"
$ws.Range("E35").Value = "  +0.25%  "

$ws.Range("D36").Value = "'0.7027"
$ws.Range("E36").Value = "  -0.18%  "

$ws.Range("D37").Value = "'2.719"
$ws.Range("E37").Value = "  +0.04%  "

$ws.Range("E38").Value = "  -0.96%  "

$ws.Range("D39").Value = "'2.763"
$ws.Range("E39").Value = "  -0.29%  "

$ws.Range("D40").Value = "'6.463"
$ws.Range("E40").Value = "  +2.35%  "

$ws.Range("D41").Value = "'70.54"
$ws.Range("E41").Value = "  -3.65%  "

$ws.Range("E42").Value = "  -1.32%  "

$ws.Range("D43").Value = "'0.8458"
$ws.Range("E43").Value = "  +0.98%  "

$ws.Range("D44").Value = "'0.4183"
$ws.Range("E44").Value = "  +0.32%  "

$ws.Range("D45").Value = "'0.9996"
$ws.Range("E45").Value = "  -0.09%  "

$ws.Range("D46").Value = "'102.90"
$ws.Range("E46").Value = "  +0.07%  "

$ws.Range("D47").Value = "'9.398"
$ws.Range("E47").Value = "  +2.64%  "

$ws.Range("D48").Value = "'7.149"
$ws.Range("E48").Value = "  +0.40%  "

$ws.Range("D49").Value = "'932.35"
$ws.Range("E49").Value = "  -4.86%  "

$ws.Range("D50").Value = "'35.30"
$ws.Range("E50").Value = "  +2.23%  "

$ws.Range("B51").Value = "Decentraland"
$ws.Range("C51").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D51").Value = "'0.3857"
$ws.Range("E51").Value = "  +0.70%  "
